$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "1.00", "269.65") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric values
# and mangles the display text (e.g. "1.00" -> 1, "269.65" -> 269.64999999999998).
$textCells = @(
    'D5',
    'D6',
    'D8',
    'D10',
    'D11',
    'D12',
    'D15',
    'D17',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D26',
    'D27',
    'D29',
    'D31',
    'D32',
    'D33',
    'D34',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all the updated coin values (price, volume %, and the occasional re-ranked
# coin name/link/price/volume swap) cell by cell, matching the refreshed feed.
$ws.Range('D2').Value = '43.676.31'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.230.89'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '269.65'
$ws.Range('E5').Value = '  +3.80%  '
$ws.Range('D6').Value = '92.68'
$ws.Range('E6').Value = '  +12.05%  '
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').Value = '46.66'
$ws.Range('E10').Value = '  +5.82%  '
$ws.Range('D11').Value = '0.0926'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = '8.25'
$ws.Range('E12').Value = '  +16.97%  '
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.567.17'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.08'
$ws.Range('E15').Value = '  +3.30%  '
$ws.Range('D16').Value = '2.218.24'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  +2.64%  '
$ws.Range('D18').Value = '43.654.81'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '6.01'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '70.42'
$ws.Range('E21').Value = '  -1.44%  '
$ws.Range('D22').Value = '2.34'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').Value = '233.29'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '8.99'
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.50'
$ws.Range('E26').Value = '  +11.34%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '11.30'
$ws.Range('E27').Value = '  +4.83%  '
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('D29').Value = '40.26'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').Value = '172.78'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '0.0921'
$ws.Range('E32').Value = '  +4.32%  '
$ws.Range('D33').Value = '20.82'
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('D34').Value = '5.48'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  -3.64%  '
$ws.Range('D37').Value = '0.0351'
$ws.Range('E37').Value = '  -2.92%  '
$ws.Range('D38').Value = '4.32'
$ws.Range('E38').Value = '  -3.76%  '
$ws.Range('D39').Value = '3.58'
$ws.Range('E39').Value = '  +20.70%  '
$ws.Range('D40').Value = '12.51'
$ws.Range('E40').Value = '  -6.84%  '
$ws.Range('E41').Value = '  +2.85%  '
$ws.Range('E42').Value = '  +8.85%  '
$ws.Range('D43').Value = '63.42'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '5.33'
$ws.Range('E44').Value = '  -3.99%  '
$ws.Range('D45').Value = '0.0989'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = '8.39'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').Value = '100.54'
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('E49').Value = '  +2.56%  '
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').Value = '2.453.67'
$ws.Range('E51').Value = '  +0.62%  '

# Drop the temporary Text number-format back to the workbook default style so the
# saved cells match the original (unstyled) data cells exactly.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
